$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.026.97"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "'3.896.62"
$ws.Range("E3").Value = "  +3.11%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'468.93"
$ws.Range("E5").Value = "  +9.62%  "
$ws.Range("D6").Value = "'144.13"
$ws.Range("E6").Value = "  +4.50%  "
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +8.03%  "
$ws.Range("E11").Value = "  +7.82%  "
$ws.Range("D12").Value = "'42.90"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'4.512.53"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'10.37"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'15.08"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'3.913.02"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").Value = "'67.199.61"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "'430.03"
$ws.Range("E21").Value = "  +5.80%  "
$ws.Range("D22").Value = "'14.74"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  +3.74%  "
$ws.Range("D24").Value = "'88.51"
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("D25").Value = "'38.64"
$ws.Range("E25").Value = "  +5.57%  "
$ws.Range("D26").Value = "'3.55"
$ws.Range("E26").Value = "  +8.02%  "
$ws.Range("E27").Value = "  +6.33%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").Value = "'9.55"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").Value = "'728.52"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "'13.70"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'43.00"
$ws.Range("E34").Value = "  +5.54%  "
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "'57.42"
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'5.40"
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("E39").Value = "  +13.83%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0476"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "'3.18"
$ws.Range("E41").Value = "  +9.29%  "
$ws.Range("D42").Value = "'0.337"
$ws.Range("E42").Value = "  +4.78%  "
$ws.Range("D43").Value = "'0.141"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +5.84%  "
$ws.Range("D46").Value = "'2.51"
$ws.Range("E46").Value = "  -6.01%  "
$ws.Range("E47").Value = "  +5.15%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  +3.22%  "
$ws.Range("D51").Value = "'143.92"
$ws.Range("E51").Value = "  +1.21%  "
